# The post that used to sit at row 287 ("「人間は一人ひとり、特別な才がある。...」")
# was removed from the source post list. Delete that entire row; every row
# below it (288..333) shifts up by one, and the sheet's used range shrinks
# from A1:C333 to A1:C332 - exactly matching the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(287).Delete()
